$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New expense entry for 08-07-2024 (SPA, 85).
# Force the new row to Text format first so values like "08-07-2024"
# and "85" are stored as literal text (matching the rest of the sheet)
# instead of being auto-converted to a date serial / number, then drop
# the temporary format again so no extra style is left behind on the row.
$newRow = $ws.Range("A18:D18")
$newRow.NumberFormat = "@"

$ws.Range("A18").Value = "08-07-2024"
$ws.Range("B18").Value = "Expense"
$ws.Range("C18").Value = "85"
$ws.Range("D18").Value = "SPA"

$newRow.ClearFormats()
